$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 180.33333
$ws.Range("I5").Value = 189.875
$ws.Range("K5").Value = 189.875
$ws.Range("M5").Value = -74.875
$ws.Range("H11").Value = 28698.535
$ws.Range("I11").Value = 28698.535
$ws.Range("K11").Value = 28698.535
$ws.Range("M11").Value = -28558.535
$ws.Range("H28").Value = 2502
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 2502
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 2502
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -3472
$ws.Range("H33").Value = 4545.8335
$ws.Range("I33").Value = 7099.067
$ws.Range("J33").Value = 290.44446
$ws.Range("K33").Value = 7099.067
$ws.Range("L33").Value = 290.44446
$ws.Range("M33").Value = -6870.067
$ws.Range("N33").Value = -748.4444599999999
$ws.Range("H51").Value = 8300.1
$ws.Range("I51").Value = 7002.5
$ws.Range("J51").Value = 9597.700000000001
$ws.Range("K51").Value = 7002.5
$ws.Range("L51").Value = 9597.700000000001
$ws.Range("M51").Value = -6518.5
$ws.Range("N51").Value = -10565.7
$ws.Range("H55").Value = 174.75
$ws.Range("I55").Value = 133
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 133
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = 81
$ws.Range("N55").Value = -728
$ws.Range("H64").Value = 3900
$ws.Range("H67").Value = 3900
$ws.Range("H86").Value = 3180.2
$ws.Range("I86").Value = 3198.5
$ws.Range("K86").Value = 3198.5
$ws.Range("M86").Value = -2075.5
$ws.Range("H89").Value = 3180.2
$ws.Range("I89").Value = 3198.5
$ws.Range("K89").Value = 15992.5
$ws.Range("M89").Value = -10376.5
$ws.Range("H98").Value = 1097
$ws.Range("I98").Value = 1139.1111
$ws.Range("K98").Value = 1139.1111
$ws.Range("M98").Value = 358.8888999999999
$ws.Range("H115").Value = 1136.1111
$ws.Range("I115").Value = 821
$ws.Range("J115").Value = 1766.3334
$ws.Range("K115").Value = 2463
$ws.Range("L115").Value = 5299.0002
$ws.Range("M115").Value = -896
$ws.Range("N115").Value = -8433.0002
$ws.Range("H122").Value = 1097
$ws.Range("I122").Value = 1139.1111
$ws.Range("K122").Value = 3417.3333
$ws.Range("M122").Value = -967.3333000000002
$ws.Range("H132").Value = 2208.4614
$ws.Range("I132").Value = 2245.9092
$ws.Range("K132").Value = 6737.7276
$ws.Range("M132").Value = -4207.7276
$ws.Range("H137").Value = 11129.934
$ws.Range("I137").Value = 1996.6
$ws.Range("J137").Value = 29396.6
$ws.Range("K137").Value = 5989.799999999999
$ws.Range("L137").Value = 88189.79999999999
$ws.Range("M137").Value = -3439.799999999999
$ws.Range("N137").Value = -93289.79999999999
$ws.Range("H138").Value = 4173.8076
$ws.Range("I138").Value = 3280.9167
$ws.Range("J138").Value = 4939.143
$ws.Range("K138").Value = 9842.750100000001
$ws.Range("L138").Value = 14817.429
$ws.Range("M138").Value = -4702.750100000001
$ws.Range("N138").Value = -25097.429
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 442.5
$ws.Range("I5").Value = 131
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 131
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -19
$ws.Range("N5").Value = -2224
$ws.Range("H31").Value = 6172.8
$ws.Range("I31").Value = 6172.8
$ws.Range("K31").Value = 6172.8
$ws.Range("M31").Value = -5878.8
$ws.Range("H32").Value = 6631.5293
$ws.Range("I32").Value = 5514
$ws.Range("K32").Value = 5514
$ws.Range("M32").Value = -5227
$ws.Range("H76").Value = 33810.5
$ws.Range("J76").Value = 33810.5
$ws.Range("L76").Value = 33810.5
$ws.Range("N76").Value = -34486.5
$ws.Range("H79").Value = 33810.5
$ws.Range("J79").Value = 33810.5
$ws.Range("L79").Value = 33810.5
$ws.Range("N79").Value = -36150.5
$ws.Range("H97").Value = 397.9655
$ws.Range("I97").Value = 363.07693
$ws.Range("J97").Value = 700.3333
$ws.Range("K97").Value = 363.07693
$ws.Range("L97").Value = 700.3333
$ws.Range("M97").Value = 132.92307
$ws.Range("N97").Value = -1692.3333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 442.5
$ws.Range("I4").Value = 131
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 131
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -16
$ws.Range("N4").Value = -2230
$ws.Range("H20").Value = 2781.1538
$ws.Range("I20").Value = 837.4286
$ws.Range("K20").Value = 837.4286
$ws.Range("M20").Value = -590.4286
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H105").Value = 2143.7173
$ws.Range("I105").Value = 1966.1786
$ws.Range("J105").Value = 2419.889
$ws.Range("K105").Value = 1966.1786
$ws.Range("L105").Value = 2419.889
$ws.Range("M105").Value = -219.1786
$ws.Range("N105").Value = -5913.889
$ws.Range("H134").Value = 4247.0435
$ws.Range("I134").Value = 4147.4116
$ws.Range("J134").Value = 4529.3335
$ws.Range("K134").Value = 12442.2348
$ws.Range("L134").Value = 13588.0005
$ws.Range("M134").Value = -9907.234800000002
$ws.Range("N134").Value = -18658.0005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2000
$ws.Range("J11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("N11").Value = -2280
$ws.Range("H25").Value = 15000
$ws.Range("H31").Value = 2400.6333
$ws.Range("I31").Value = 1519.8572
$ws.Range("J31").Value = 3171.3125
$ws.Range("K31").Value = 1519.8572
$ws.Range("L31").Value = 3171.3125
$ws.Range("M31").Value = -1224.8572
$ws.Range("N31").Value = -3761.3125
$ws.Range("H34").Value = 2400.6333
$ws.Range("I34").Value = 1519.8572
$ws.Range("J34").Value = 3171.3125
$ws.Range("K34").Value = 1519.8572
$ws.Range("L34").Value = 3171.3125
$ws.Range("M34").Value = -1317.8572
$ws.Range("N34").Value = -3575.3125
$ws.Range("H50").Value = 29648.572
$ws.Range("J50").Value = 29999.846
$ws.Range("L50").Value = 29999.846
$ws.Range("N50").Value = -31249.846
$ws.Range("H58").Value = 2471.9285
$ws.Range("I58").Value = 2213.125
$ws.Range("K58").Value = 2213.125
$ws.Range("M58").Value = -2010.125
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H70").Value = 99885.39999999999
$ws.Range("J70").Value = 99885.39999999999
$ws.Range("L70").Value = 99885.39999999999
$ws.Range("N70").Value = -100515.4
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H73").Value = 99885.39999999999
$ws.Range("J73").Value = 99885.39999999999
$ws.Range("L73").Value = 99885.39999999999
$ws.Range("N73").Value = -102069.4
$ws.Range("H75").Value = 39990
$ws.Range("J75").Value = 39990
$ws.Range("L75").Value = 39990
$ws.Range("N75").Value = -41986
$ws.Range("H78").Value = 39990
$ws.Range("J78").Value = 39990
$ws.Range("L78").Value = 119970
$ws.Range("N78").Value = -129954
$ws.Range("H86").Value = 4923.643
$ws.Range("I86").Value = 4912.4546
$ws.Range("J86").Value = 4964.6665
$ws.Range("K86").Value = 4912.4546
$ws.Range("L86").Value = 4964.6665
$ws.Range("M86").Value = -3789.4546
$ws.Range("N86").Value = -7210.6665
$ws.Range("H89").Value = 4923.643
$ws.Range("I89").Value = 4912.4546
$ws.Range("J89").Value = 4964.6665
$ws.Range("K89").Value = 24562.273
$ws.Range("L89").Value = 24823.3325
$ws.Range("M89").Value = -18946.273
$ws.Range("N89").Value = -36055.3325
$ws.Range("H122").Value = 2274.75
$ws.Range("I122").Value = 2433
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 7299
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -4849
$ws.Range("N122").Value = -10300
$ws.Range("H132").Value = 5151
$ws.Range("I132").Value = 4655.8
$ws.Range("J132").Value = 9277.666999999999
$ws.Range("K132").Value = 13967.4
$ws.Range("L132").Value = 27833.001
$ws.Range("M132").Value = -11437.4
$ws.Range("N132").Value = -32893.001
$ws.Range("H134").Value = 3866.1428
$ws.Range("I134").Value = 4158.294
$ws.Range("J134").Value = 2624.5
$ws.Range("K134").Value = 12474.882
$ws.Range("L134").Value = 7873.5
$ws.Range("M134").Value = -9939.882
$ws.Range("N134").Value = -12943.5
$ws.Range("H136").Value = 2471.9285
$ws.Range("I136").Value = 2213.125
$ws.Range("K136").Value = 6639.375
$ws.Range("M136").Value = -4089.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1765.8334
$ws.Range("J17").Value = 2449
$ws.Range("L17").Value = 7347
$ws.Range("N17").Value = -7685
$ws.Range("H22").Value = 1374.6666
$ws.Range("I22").Value = 395.5
$ws.Range("J22").Value = 3333
$ws.Range("K22").Value = 1186.5
$ws.Range("L22").Value = 9999
$ws.Range("M22").Value = -1017.5
$ws.Range("N22").Value = -10337
$ws.Range("H27").Value = 1374.6666
$ws.Range("I27").Value = 395.5
$ws.Range("J27").Value = 3333
$ws.Range("K27").Value = 1186.5
$ws.Range("L27").Value = 9999
$ws.Range("M27").Value = -1084.5
$ws.Range("N27").Value = -10203
$ws.Range("H29").Value = 2814.1667
$ws.Range("I29").Value = 1999
$ws.Range("K29").Value = 5997
$ws.Range("M29").Value = -5720
$ws.Range("H32").Value = 900
$ws.Range("I32").Value = 900
$ws.Range("K32").Value = 2700
$ws.Range("M32").Value = -2417
$ws.Range("H33").Value = 132.22223
$ws.Range("I33").Value = 240
$ws.Range("J33").Value = 118.75
$ws.Range("K33").Value = 1440
$ws.Range("L33").Value = 712.5
$ws.Range("M33").Value = -1157
$ws.Range("N33").Value = -1278.5
$ws.Range("H34").Value = 3787.8125
$ws.Range("I34").Value = 313
$ws.Range("J34").Value = 7262.625
$ws.Range("K34").Value = 939
$ws.Range("L34").Value = 21787.875
$ws.Range("M34").Value = -855
$ws.Range("N34").Value = -21955.875
$ws.Range("H68").Value = 1639.7142
$ws.Range("I68").Value = 1741
$ws.Range("K68").Value = 5223
$ws.Range("M68").Value = -4412
$ws.Range("H71").Value = 1639.7142
$ws.Range("I71").Value = 1741
$ws.Range("K71").Value = 15669
$ws.Range("M71").Value = -11613
$ws.Range("H131").Value = 1480.1428
$ws.Range("I131").Value = 556.8333
$ws.Range("J131").Value = 1608.9767
$ws.Range("K131").Value = 1670.4999
$ws.Range("L131").Value = 4826.9301
$ws.Range("M131").Value = 3369.5001
$ws.Range("N131").Value = -14906.9301
$ws.Range("H138").Value = 1788
$ws.Range("I138").Value = 1788
$ws.Range("K138").Value = 5364
$ws.Range("M138").Value = -224
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 40999.5
$ws.Range("J63").Value = 40999.5
$ws.Range("L63").Value = 40999.5
$ws.Range("N63").Value = -42371.5
$ws.Range("H66").Value = 40999.5
$ws.Range("J66").Value = 40999.5
$ws.Range("L66").Value = 122998.5
$ws.Range("N66").Value = -129862.5
$ws.Range("H70").Value = 7236
$ws.Range("I70").Value = 5942.75
$ws.Range("K70").Value = 5942.75
$ws.Range("M70").Value = -5672.75
$ws.Range("H73").Value = 7236
$ws.Range("I73").Value = 5942.75
$ws.Range("K73").Value = 5942.75
$ws.Range("M73").Value = -5006.75
$ws.Range("H113").Value = 252626.5
$ws.Range("J113").Value = 3013
$ws.Range("L113").Value = 3013
$ws.Range("N113").Value = -7353
$ws.Range("H122").Value = 1886.3478
$ws.Range("I122").Value = 1195.8422
$ws.Range("K122").Value = 3587.5266
$ws.Range("M122").Value = -1137.5266
$ws.Range("H126").Value = 3761.682
$ws.Range("I126").Value = 3398
$ws.Range("J126").Value = 3868.647
$ws.Range("K126").Value = 10194
$ws.Range("L126").Value = 11605.941
$ws.Range("M126").Value = -7724
$ws.Range("N126").Value = -16545.941
$ws.Range("H132").Value = 3908.7058
$ws.Range("I132").Value = 3571.8462
$ws.Range("K132").Value = 10715.5386
$ws.Range("M132").Value = -8185.5386
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 15000
$ws.Range("J14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("N14").Value = -15344
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H43").Value = 462327.1
$ws.Range("J43").Value = 505115.38
$ws.Range("L43").Value = 505115.38
$ws.Range("N43").Value = -505501.38
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 3271041.8
$ws.Range("I132").Value = 3586884.8
$ws.Range("K132").Value = 10760654.4
$ws.Range("M132").Value = -10758124.4
$ws.Range("H136").Value = 11113775
$ws.Range("I136").Value = 12348528
$ws.Range("K136").Value = 37045584
$ws.Range("M136").Value = -37043034
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H110").Value = 64990
$ws.Range("J110").Value = 64990
$ws.Range("L110").Value = 64990
$ws.Range("N110").Value = -73170
$ws.Range("H126").Value = 4629.359
$ws.Range("I126").Value = 4549.2593
$ws.Range("K126").Value = 13647.7779
$ws.Range("M126").Value = -11177.7779
$ws.Range("H132").Value = 12267.096
$ws.Range("I132").Value = 8367.166999999999
$ws.Range("K132").Value = 25101.501
$ws.Range("M132").Value = -22571.501
$ws.Range("H136").Value = 4546738
$ws.Range("I136").Value = 4763240
$ws.Range("J136").Value = 200
$ws.Range("K136").Value = 14289720
$ws.Range("L136").Value = 600
$ws.Range("M136").Value = -14287170
$ws.Range("N136").Value = -5700
